# The edit rotates the three species-observation records currently stored
# in rows 8, 9 and 10 of the "Artfynd" sheet:
#   new row 8  <- old row 10 (Trådticka / Climacocystis borealis)
#   new row 9  <- old row 8  (Rosenticka / Rhodofomes roseus)
#   new row 10 <- old row 9  (Garnlav / Alectoria sarmentosa)
# Only the per-record columns are touched; columns that are identical for
# all three rows (C, P, S, T, U, V, W, Y, AA, AD, AE, AG, AJ, AK, AT, AW,
# AX, AY, ...) are left untouched since rotating them would be a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8  (becomes the old row-10 record: Trådticka) ---------------------
$ws.Range("A8").Value = 112038082
$ws.Range("B8").Value = 90235
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 3298
$ws.Range("F8").Value = 'Trådticka'
$ws.Range("G8").Value = 'Climacocystis borealis'
$ws.Range("H8").Value = '(Fr.) Kotl. & Pouzar'
# Row 10 had no Antal/Enhet ("I"/"J") values, so row 8 loses them too.
$ws.Range("I8").Value = "'"
$ws.Range("J8").ClearContents()
$ws.Range("Q8").Value = 515925
$ws.Range("R8").Value = 7184319
$ws.Range("Z8").Value = '13:22'
$ws.Range("AB8").Value = '13:22'
$ws.Range("AM8").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO8").Value = 'Standing dead tree/snags # Picea abies'

# --- Row 9  (becomes the old row-8 record: Rosenticka) ---------------------
$ws.Range("A9").Value = 112038473
$ws.Range("B9").Value = 89834
$ws.Range("E9").Value = 658
$ws.Range("F9").Value = 'Rosenticka'
$ws.Range("G9").Value = 'Rhodofomes roseus'
$ws.Range("H9").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
# Leading "'" forces the numeric-looking Antal value to stay text, matching
# how this workbook always stores column I.
$ws.Range("I9").Value = "'4"
$ws.Range("J9").Value = 'fruktkroppar'
$ws.Range("Q9").Value = 516057
$ws.Range("R9").Value = 7184320
$ws.Range("Z9").Value = '13:34'
$ws.Range("AB9").Value = '13:34'
$ws.Range("AH9").Value = 'Blåbärsgranskog'
$ws.Range("AM9").Value = 'Liggande död trädstam, markontakt'
$ws.Range("AO9").Value = 'Horizontal, dead with ground contact # Picea abies'

# --- Row 10 (becomes the old row-9 record: Garnlav) -------------------------
$ws.Range("A10").Value = 112037684
$ws.Range("B10").Value = 77650
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = 'Garnlav'
$ws.Range("G10").Value = 'Alectoria sarmentosa'
$ws.Range("H10").Value = '(Ach.) Ach.'
$ws.Range("Q10").Value = 515886
$ws.Range("R10").Value = 7184226
$ws.Range("Z10").Value = '12:08'
$ws.Range("AB10").Value = '12:08'
$ws.Range("AH10").Value = 'Blåbärsbarrskog'
